# daily auto push: 2026-01-21 03:50 UTC
# A new sample ("2026/01/21" 水 11 201) was recorded for the day, inserted
# as a new row 684, pushing every following row (old 684..725) down by one
# (new 685..726). The sheet's used range grows from D725 to D726.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 684; Excel shifts rows 684:725 down to 685:726,
# carrying their values/formatting with them intact.
$ws.Rows("684:684").Insert()

# Force column A to be stored as text (not auto-parsed into a date serial),
# matching the rest of the "日付" column which holds literal "yyyy/mm/dd" strings.
$ws.Range("A684").NumberFormat = "@"
$ws.Range("A684").Value = "2026/01/21"
# Drop the temporary text number-format again so the new row ends up with
# the same (default/general) cell style as every other data row.
$ws.Range("A684").ClearFormats()

$ws.Range("B684").Value = "水"
$ws.Range("C684").Value = 11
$ws.Range("D684").Value = 201

Write-Host "done"
